$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (NumberError): add 6 new rows (62-67), copying formatting from row 61 ---
$ws1.Range("A61:G61").Copy()
$ws1.Range("A62:G67").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Sheet1: populate new rows 62-67 ---
$ws1.Range("A62").Value = "Pass"
$ws1.Range("B62").Value = "Tue Feb 11 20:06:41 EST 2025"
$ws1.Range("C62").Value = "Y"
$ws1.Range("D62").Value = "Existing Liability with Notice/Invoice Number"
$ws1.Range("E62").Value = "Digital Advertising Gross Revenues"
$ws1.Range("F62").Value = "123456789"
$ws1.Range("G62").Value = "Notice Number or Invoice Number must be 10, 11, or 13 digits in length, with a value greater than zero"

$ws1.Range("A63").Value = "Pass"
$ws1.Range("B63").Value = "Tue Feb 11 20:06:50 EST 2025"
$ws1.Range("C63").Value = "Y"
$ws1.Range("D63").Value = "Existing Liability with Notice/Invoice Number"
$ws1.Range("E63").Value = "Digital Advertising Gross Revenues"
$ws1.Range("F63").Value = "123456789012"
$ws1.Range("G63").Value = "Notice Number or Invoice Number must be 10, 11, or 13 digits in length, with a value greater than zero"

$ws1.Range("A64").Value = "Pass"
$ws1.Range("B64").Value = "Tue Feb 11 20:07:00 EST 2025"
$ws1.Range("C64").Value = "Y"
$ws1.Range("D64").Value = "Existing Liability with Notice/Invoice Number"
$ws1.Range("E64").Value = "Digital Advertising Gross Revenues"
$ws1.Range("F64").Value = "0000000000"
$ws1.Range("G64").Value = "Notice Number or Invoice Number must be 10, 11, or 13 digits in length, with a value greater than zero"

$ws1.Range("A65").Value = "Pass"
$ws1.Range("B65").Value = "Tue Feb 11 20:07:09 EST 2025"
$ws1.Range("C65").Value = "Y"
$ws1.Range("D65").Value = "Existing Liability with Notice/Invoice Number"
$ws1.Range("E65").Value = "PTE Composite"
$ws1.Range("F65").Value = "123456789"
$ws1.Range("G65").Value = "Notice Number or Invoice Number must be 10, 11, or 13 digits in length, with a value greater than zero"

$ws1.Range("A66").Value = "Pass"
$ws1.Range("B66").Value = "Tue Feb 11 20:07:18 EST 2025"
$ws1.Range("C66").Value = "Y"
$ws1.Range("D66").Value = "Existing Liability with Notice/Invoice Number"
$ws1.Range("E66").Value = "PTE Composite"
$ws1.Range("F66").Value = "123456789012"
$ws1.Range("G66").Value = "Notice Number or Invoice Number must be 10, 11, or 13 digits in length, with a value greater than zero"

$ws1.Range("A67").Value = "Pass"
$ws1.Range("B67").Value = "Tue Feb 11 20:07:27 EST 2025"
$ws1.Range("C67").Value = "Y"
$ws1.Range("D67").Value = "Existing Liability with Notice/Invoice Number"
$ws1.Range("E67").Value = "PTE Composite"
$ws1.Range("F67").Value = "0000000000"
$ws1.Range("G67").Value = "Notice Number or Invoice Number must be 10, 11, or 13 digits in length, with a value greater than zero"

# --- Sheet1: update Date (column B) for rows 2-61 with new timestamps ---
$ws1.Range("B2").Value = "Tue Feb 11 19:57:26 EST 2025"
$ws1.Range("B3").Value = "Tue Feb 11 19:57:35 EST 2025"
$ws1.Range("B4").Value = "Tue Feb 11 19:57:45 EST 2025"
$ws1.Range("B5").Value = "Tue Feb 11 19:57:54 EST 2025"
$ws1.Range("B6").Value = "Tue Feb 11 19:58:04 EST 2025"
$ws1.Range("B7").Value = "Tue Feb 11 19:58:13 EST 2025"
$ws1.Range("B8").Value = "Tue Feb 11 19:58:22 EST 2025"
$ws1.Range("B9").Value = "Tue Feb 11 19:58:31 EST 2025"
$ws1.Range("B10").Value = "Tue Feb 11 19:58:40 EST 2025"
$ws1.Range("B11").Value = "Tue Feb 11 19:58:49 EST 2025"
$ws1.Range("B12").Value = "Tue Feb 11 19:58:58 EST 2025"
$ws1.Range("B13").Value = "Tue Feb 11 19:59:08 EST 2025"
$ws1.Range("B14").Value = "Tue Feb 11 19:59:17 EST 2025"
$ws1.Range("B15").Value = "Tue Feb 11 19:59:26 EST 2025"
$ws1.Range("B16").Value = "Tue Feb 11 19:59:35 EST 2025"
$ws1.Range("B17").Value = "Tue Feb 11 19:59:44 EST 2025"
$ws1.Range("B18").Value = "Tue Feb 11 19:59:53 EST 2025"
$ws1.Range("B19").Value = "Tue Feb 11 20:00:03 EST 2025"
$ws1.Range("B20").Value = "Tue Feb 11 20:00:12 EST 2025"
$ws1.Range("B21").Value = "Tue Feb 11 20:00:21 EST 2025"
$ws1.Range("B22").Value = "Tue Feb 11 20:00:31 EST 2025"
$ws1.Range("B23").Value = "Tue Feb 11 20:00:40 EST 2025"
$ws1.Range("B24").Value = "Tue Feb 11 20:00:49 EST 2025"
$ws1.Range("B25").Value = "Tue Feb 11 20:00:59 EST 2025"
$ws1.Range("B26").Value = "Tue Feb 11 20:01:08 EST 2025"
$ws1.Range("B27").Value = "Tue Feb 11 20:01:17 EST 2025"
$ws1.Range("B28").Value = "Tue Feb 11 20:01:26 EST 2025"
$ws1.Range("B29").Value = "Tue Feb 11 20:01:35 EST 2025"
$ws1.Range("B30").Value = "Tue Feb 11 20:01:44 EST 2025"
$ws1.Range("B31").Value = "Tue Feb 11 20:01:54 EST 2025"
$ws1.Range("B32").Value = "Tue Feb 11 20:02:03 EST 2025"
$ws1.Range("B33").Value = "Tue Feb 11 20:02:12 EST 2025"
$ws1.Range("B34").Value = "Tue Feb 11 20:02:22 EST 2025"
$ws1.Range("B35").Value = "Tue Feb 11 20:02:31 EST 2025"
$ws1.Range("B36").Value = "Tue Feb 11 20:02:40 EST 2025"
$ws1.Range("B37").Value = "Tue Feb 11 20:02:49 EST 2025"
$ws1.Range("B38").Value = "Tue Feb 11 20:02:59 EST 2025"
$ws1.Range("B39").Value = "Tue Feb 11 20:03:08 EST 2025"
$ws1.Range("B40").Value = "Tue Feb 11 20:03:17 EST 2025"
$ws1.Range("B41").Value = "Tue Feb 11 20:03:26 EST 2025"
$ws1.Range("B42").Value = "Tue Feb 11 20:03:36 EST 2025"
$ws1.Range("B43").Value = "Tue Feb 11 20:03:45 EST 2025"
$ws1.Range("B44").Value = "Tue Feb 11 20:03:54 EST 2025"
$ws1.Range("B45").Value = "Tue Feb 11 20:04:03 EST 2025"
$ws1.Range("B46").Value = "Tue Feb 11 20:04:14 EST 2025"
$ws1.Range("B47").Value = "Tue Feb 11 20:04:23 EST 2025"
$ws1.Range("B48").Value = "Tue Feb 11 20:04:32 EST 2025"
$ws1.Range("B49").Value = "Tue Feb 11 20:04:42 EST 2025"
$ws1.Range("B50").Value = "Tue Feb 11 20:04:51 EST 2025"
$ws1.Range("B51").Value = "Tue Feb 11 20:05:00 EST 2025"
$ws1.Range("B52").Value = "Tue Feb 11 20:05:09 EST 2025"
$ws1.Range("B53").Value = "Tue Feb 11 20:05:18 EST 2025"
$ws1.Range("B54").Value = "Tue Feb 11 20:05:28 EST 2025"
$ws1.Range("B55").Value = "Tue Feb 11 20:05:37 EST 2025"
$ws1.Range("B56").Value = "Tue Feb 11 20:05:46 EST 2025"
$ws1.Range("B57").Value = "Tue Feb 11 20:05:55 EST 2025"
$ws1.Range("B58").Value = "Tue Feb 11 20:06:04 EST 2025"
$ws1.Range("B59").Value = "Tue Feb 11 20:06:13 EST 2025"
$ws1.Range("B60").Value = "Tue Feb 11 20:06:23 EST 2025"
$ws1.Range("B61").Value = "Tue Feb 11 20:06:32 EST 2025"

# --- View state: selections / active sheet ---
[void]$ws2.Range("E25").Select()
$ws1.Activate()
[void]$ws1.Range("C59:C67").Select()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
